# Revert to 2.1.1 files
# The OCCF workbook's "About" sheet references a 2019-dollar CPI conversion
# factor; this change reverts those labels/values back to the 2018-dollar
# baseline used in v2.1.1 of the model.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Large Output Currency Unit scale label: "billion 2019 dollars" -> "billion 2018 dollars"
$ws.Range("A18").Value = "billion 2018 dollars"

# Medium Output Currency Unit scale label: "million 2019 dollars" -> "million 2018 dollars"
$ws.Range("A21").Value = "million 2018 dollars"

# Conversion-factor caption: "2019 dollars per 2012 dollar" -> "2018 dollars per 2012 dollar"
$ws.Range("B26").Value = "2018 dollars per 2012 dollar"

# Explanatory note referencing the same ratio
$ws.Range("B29").Value = "which in this case is ""2012 dollars per 2018 dollar."""

# The conversion factor itself (2012 dollars per 2018 dollar, from cpi.xlsx)
$ws.Range("A26").Value = 0.9143273584567535
